$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2023-01-04"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "#CLR99L9L"
$ws.Range("D2").Value = 38197
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = "Mythic I"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = "AMBER, ASH, BARLEY, BEA, BELLE, BO, BONNIE, BROCK, BYRON, CARL, COLETTE, CROW, ... "
$ws.Range("A3").Value = "player"
$ws.Range("B3").Value = "team"
$ws.Range("C3").Value = "tag"
$ws.Range("D3").Value = "trophies"
$ws.Range("E3").Value = "pl_score"
$ws.Range("F3").Value = "pl_rank"
$ws.Range("G3").Value = "level_9s"
$ws.Range("H3").Value = "level_10s"
$ws.Range("I3").Value = "level_11s"
$ws.Range("J3").Value = "brawlers_11"
$ws.Range("A4").Value = "PFCXF"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "#CLR99L9L"
$ws.Range("D4").Value = 38197
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Mythic I"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 24
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = "AMBER, ASH, BARLEY, BEA, BELLE, BO, BONNIE, BROCK, BYRON, CARL, COLETTE, CROW, ... "
$ws.Range("A5").Value = "Sean747❤️Alexa"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "#QYJGL9QU"
$ws.Range("D5").Value = 6558
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Mythic I"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 21
$ws.Range("J5").Value = "BELLE, BIBI, BROCK, BULL, CARL, COLT, CROW, EL PRIMO, JESSIE, MAX, MORTIS, MR. P, ... "
$ws.Range("A6").Value = "VS | dragoh"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "#GC0LRQ0V"
$ws.Range("D6").Value = 28808
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 27
$ws.Range("I6").Value = 26
$ws.Range("J6").Value = "8-BIT, AMBER, ASH, BARLEY, BELLE, BROCK, BYRON, CARL, DARRYL, GALE, GENE, GRIFF, ... "
$ws.Range("A7").Value = "BennyBoy_"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "#2PR80P8CU"
$ws.Range("D7").Value = 49867
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 19
$ws.Range("I7").Value = 32
$ws.Range("J7").Value = "AMBER, BEA, BO, BUZZ, BYRON, CARL, CHESTER, COLETTE, CROW, DARRYL, DYNAMIKE, EMZ, ... "
$ws.Range("A8").Value = "Tribe|LHC"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "#J2RLUJP2"
$ws.Range("D8").Value = 48129
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = "Mythic III"
$ws.Range("G8").Value = 14
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 40
$ws.Range("J8").Value = "8-BIT, AMBER, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BUSTER, BUZZ, BYRON, CARL, ... "
$ws.Range("A9").Value = "Zachary ϟ"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "#RQ8RRV0Y"
$ws.Range("D9").Value = 33101
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "Diamond I"
$ws.Range("G9").Value = 17
$ws.Range("H9").Value = 35
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = "BYRON, MAX, MORTIS, RICO, SPIKE, STU, TARA"
$ws.Range("A10").Value = "David"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "#9UPQVRQ"
$ws.Range("D10").Value = 36106
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = "Silver III"
$ws.Range("G10").Value = 37
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = "BROCK, EMZ, EVE, GALE, ROSA, SPROUT, TICK"
$ws.Range("A11").Value = "Dead Beat"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "#2LJGQ02YV"
$ws.Range("D11").Value = 38678
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 24
$ws.Range("I11").Value = 32
$ws.Range("J11").Value = "AMBER, ASH, BEA, BELLE, BONNIE, BUSTER, CARL, CHESTER, COLETTE, CROW, GALE, GENE, ... "
$ws.Range("A12").Value = "Oh No! Jotaro"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "#YQ0LJU2"
$ws.Range("D12").Value = 50845
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = "Mythic I"
$ws.Range("G12").Value = 24
$ws.Range("H12").Value = 18
$ws.Range("I12").Value = 21
$ws.Range("J12").Value = "8-BIT, AMBER, BEA, BELLE, BIBI, COLETTE, EMZ, GENE, GRAY, GUS, JACKY, JANET, ... "
$ws.Range("A13").Value = "Cred"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "#RGVL0L8V"
$ws.Range("D13").Value = 42374
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 62
$ws.Range("J13").Value = "8-BIT, AMBER, ASH, BARLEY, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BULL, BUSTER, ... "
$ws.Range("A14").Value = "DepressedRekt™"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = "#VY9JJ9RQ"
$ws.Range("D14").Value = 34713
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = "Mythic I"
$ws.Range("G14").Value = 16
$ws.Range("H14").Value = 39
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = "BEA, BELLE, GENE, GUS, MAX, PIPER, SPIKE, STU"
$ws.Range("A15").Value = "哥本哈根大使"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "#9CQURRU28"
$ws.Range("D15").Value = 43562
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = "Mythic II"
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 24
$ws.Range("I15").Value = 27
$ws.Range("J15").Value = "ASH, BEA, BELLE, BIBI, BO, BROCK, BYRON, CARL, CHESTER, COLETTE, FANG, GENE, ... "
$ws.Range("A16").Value = "Blue"
$ws.Range("B16").Value = 5
$ws.Range("C16").Value = "#C9LR0R0V"
$ws.Range("D16").Value = 40869
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Mythic I"
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = "AMBER, ASH, BEA, BO, BONNIE, BUSTER, BYRON, CHESTER, COLETTE, CROW, EMZ, EVE, ... "
$ws.Range("A17").Value = "Hogglific"
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = "#UCY09URC"
$ws.Range("D17").Value = 45367
$ws.Range("E17").Value = 14
$ws.Range("F17").Value = "Mythic II"
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = "8-BIT, AMBER, ASH, BARLEY, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BULL, BUSTER, ... "
$ws.Range("A18").Value = "ZSilverZ"
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = "#2YQUPUYJ"
$ws.Range("D18").Value = 49185
$ws.Range("E18").Value = 14
$ws.Range("F18").Value = "Mythic II"
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 9
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = "8-BIT, AMBER, ASH, BEA, BELLE, BIBI, BO, BROCK, BULL, BUSTER, BUZZ, CARL, ... "
$ws.Range("A19").Value = "prasanna"
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = "#Y09QC2UJQ"
$ws.Range("D19").Value = 35142
$ws.Range("E19").Value = 11
$ws.Range("F19").Value = "Diamond II"
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 17
$ws.Range("I19").Value = 12
$ws.Range("J19").Value = "BUSTER, CHESTER, EL PRIMO, EMZ, FANG, GRIFF, LOLA, MAX, NITA, OTIS, SANDY"
$ws.Range("A20").Value = "❄|Vırtual ◇"
$ws.Range("B20").Value = 6
$ws.Range("C20").Value = "#U99LR8PR"
$ws.Range("D20").Value = 34097
$ws.Range("E20").Value = 14
$ws.Range("F20").Value = "Mythic II"
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 13
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = "ASH, BELLE, BIBI, CARL, COLETTE, EMZ, GENE, MAX, RICO, STU, SURGE"
$ws.Range("A21").Value = "冰Ferds ☃️"
$ws.Range("B21").Value = 6
$ws.Range("C21").Value = "#YPRVGV8L"
$ws.Range("D21").Value = 35749
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Gold III"
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 26
$ws.Range("I21").Value = 10
$ws.Range("J21").Value = "BELLE, BIBI, BULL, FRANK, MAX, MORTIS, RICO, SQUEAK, STU, TARA"
$ws.Range("A22").Value = "Saurav"
$ws.Range("B22").Value = 7
$ws.Range("C22").Value = "#U2Q9L2QU"
$ws.Range("D22").Value = 34843
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = ""
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = "BEA, CARL, COLT, EMZ, GENE, JACKY, MAX, MORTIS, SURGE, TARA"
$ws.Range("A23").Value = "VS∣Daniel"
$ws.Range("B23").Value = 7
$ws.Range("C23").Value = "#PCRG080"
$ws.Range("D23").Value = 49820
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Mythic I"
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 39
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = "BEA, BELLE, CROW, EMZ, FANG, GENE, LEON, MAX, MORTIS, RUFFS, SANDY, SQUEAK, ... "
$ws.Range("A24").Value = "majd✿"
$ws.Range("B24").Value = 7
$ws.Range("C24").Value = "#PY0YR2CR8"
$ws.Range("D24").Value = 40895
$ws.Range("E24").Value = 14
$ws.Range("F24").Value = "Mythic II"
$ws.Range("G24").Value = 23
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 8
$ws.Range("J24").Value = "BEA, BELLE, CARL, CHESTER, MORTIS, PENNY, PIPER, STU"
$ws.Range("A25").Value = "KaiWen"
$ws.Range("B25").Value = 8
$ws.Range("C25").Value = "#YUCGRY2V"
$ws.Range("D25").Value = 35940
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 14
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = "SHELLY"
$ws.Range("A26").Value = "SG|🔥Fire🔥"
$ws.Range("B26").Value = 8
$ws.Range("C26").Value = "#CQ8RC802"
$ws.Range("D26").Value = 44940
$ws.Range("E26").Value = 12
$ws.Range("F26").Value = "Diamond III"
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 27
$ws.Range("I26").Value = 8
$ws.Range("J26").Value = "BUZZ, CARL, CHESTER, GENE, MAX, MORTIS, SPIKE, SURGE"
$ws.Range("A27").Value = "꧁✨やͪeͧ༒ͨ 𝖊𝖗✨꧂"
$ws.Range("B27").Value = 8
$ws.Range("C27").Value = "#RQUQG8JY"
$ws.Range("D27").Value = 39681
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Mythic I"
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 18
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = "BEA, BELLE, BIBI, BO, BROCK, BUSTER, BUZZ, BYRON, CARL, CHESTER, COLETTE, CROW, ... "
$ws.Range("A28").Value = "JustBeNice"
$ws.Range("B28").Value = 9
$ws.Range("C28").Value = "#J90L8C0Y"
$ws.Range("D28").Value = 34367
$ws.Range("E28").Value = 10
$ws.Range("F28").Value = "Diamond I"
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 53
$ws.Range("J28").Value = "8-BIT, AMBER, ASH, BARLEY, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BULL, BUZZ, ... "
$ws.Range("A29").Value = "Procrastinator⏳"
$ws.Range("B29").Value = 9
$ws.Range("C29").Value = "#V2JCQQQJ"
$ws.Range("D29").Value = 37897
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 52
$ws.Range("J29").Value = "8-BIT, AMBER, ASH, BARLEY, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BUZZ, BYRON, ... "
$ws.Range("A30").Value = "menoobenelol"
$ws.Range("B30").Value = 9
$ws.Range("C30").Value = "#CP988U8Y"
$ws.Range("D30").Value = 33861
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = ""
$ws.Range("G30").Value = 28
$ws.Range("H30").Value = 3
$ws.Range("I30").Value = 29
$ws.Range("J30").Value = "ASH, BELLE, BO, BONNIE, BROCK, BUZZ, BYRON, CARL, COLETTE, COLT, CROW, DARRYL, ... "
$ws.Range("A31").Value = "XiXi"
$ws.Range("B31").Value = 10
$ws.Range("C31").Value = "#CJ2Y202J"
$ws.Range("D31").Value = 34298
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = "Diamond I"
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 11
$ws.Range("I31").Value = 32
$ws.Range("J31").Value = "BELLE, BO, BONNIE, BROCK, BUSTER, BUZZ, BYRON, CARL, CHESTER, DARRYL, EMZ, EVE, ... "
$ws.Range("A32").Value = "gkb"
$ws.Range("B32").Value = 10
$ws.Range("C32").Value = "#J9CJGRLG"
$ws.Range("D32").Value = 38145
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = "Gold I"
$ws.Range("G32").Value = 18
$ws.Range("H32").Value = 37
$ws.Range("I32").Value = 8
$ws.Range("J32").Value = "BROCK, BYRON, CROW, LEON, MAX, PENNY, PIPER, STU"
$ws.Range("A33").Value = "one"
$ws.Range("B33").Value = 10
$ws.Range("C33").Value = "#2LQ899P82"
$ws.Range("D33").Value = 41336
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = "Diamond I"
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 15
$ws.Range("I33").Value = 43
$ws.Range("J33").Value = "8-BIT, ASH, BARLEY, BEA, BELLE, BO, BONNIE, BROCK, BYRON, CARL, COLETTE, COLT, ... "
